# Append two new signup rows (4 and 5) to the "User Signups" sheet,
# matching the style of the existing alternating-row formatting
# (row 4 mirrors row 2's style, row 5 mirrors row 3's style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Test User signup
$ws.Range("A4").Value = "2025-10-26T09:48:47.915Z"
$ws.Range("B4").Value = "Test User"
$ws.Range("C4").Value = "test@example.com"
$ws.Range("D4").Value = "testuser123"
$ws.Range("E4").Value = "password123"
$ws.Range("F4").Value = "::1"
$ws.Range("G4").Value = "Active"

# Row 5: Harsh Chandrakant Mali signup
$ws.Range("A5").Value = "2025-10-26T09:49:40.649Z"
$ws.Range("B5").Value = "Harsh Chandrakant Mali"
$ws.Range("C5").Value = "jyxuta@cyclelove.cc"
$ws.Range("D5").Value = "Q23n1sdjk"
$ws.Range("E5").Value = "qwe123decdcs"
$ws.Range("F5").Value = "::1"
$ws.Range("G5").Value = "Active"

# Carry over the existing alternating row formatting (fill/border/font)
# from rows 2 and 3 onto the new rows 4 and 5 without introducing new
# cell styles.
$ws.Range("A2:G2").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)

$ws.Range("A3:G3").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
